$wb = $excel.ActiveWorkbook

# ---------- Sheet 'R2' ----------
$ws1 = $wb.Worksheets.Item("R2")
$ws1.Cells.Item(1, 1).Value = "Degree_Rede_EntreajudaLabur"
$ws1.Cells.Item(1, 2).Value = 0.4656312315740796
$ws1.Cells.Item(2, 1).Value = "OutDeg_Var.Dep_RedeControlExtAusencia"
$ws1.Cells.Item(2, 2).Value = 0.5557770026305787
$ws1.Cells.Item(3, 1).Value = "OutDeg_Var.Dep_ApoioSpecDiqCint"
$ws1.Cells.Item(3, 2).Value = 0.8847674119520766
$ws1.Cells.Item(4, 1).Value = "OutDeg_Var.Dep_ApoioSpecTubos"
$ws1.Cells.Item(4, 2).Value = 0.9813066316034472
$ws1.Cells.Item(5, 1).Value = "OutDeg_Var.Dep_RedApoiLevntDiCin"
$ws1.Cells.Item(5, 2).Value = 0.5483602824544844
$ws1.Cells.Item(6, 1).Value = "OutDeg_RedCont_DiqPriqDiqCint"
$ws1.Cells.Item(6, 2).Value = 0.8516980174315374
$ws1.Cells.Item(7, 1).Value = "OutDeg_RedeRepar_DiqPriqDiqCint"
$ws1.Cells.Item(7, 2).Value = 0.7731602918251106
$ws1.Cells.Item(8, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqCin"
$ws1.Cells.Item(8, 2).Value = 0.6902441221495704
$ws1.Cells.Item(9, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqPrq"
$ws1.Cells.Item(9, 2).Value = 0.7980698830403878
$ws1.Cells.Item(10, 1).Value = "OutDeg_ContRep2_DiqPrqCin"
$ws1.Cells.Item(10, 2).Value = 0.8259688774596836
$ws1.Cells.Item(11, 1).Value = "OutDeg_Var.Dep_RedePartilhaAgu"
$ws1.Cells.Item(11, 2).Value = 0.4132278312217763

# ---------- Sheet 'R2 Ajustado' ----------
$ws2 = $wb.Worksheets.Item("R2 Ajustado")
$ws2.Cells.Item(1, 1).Value = "Degree_Rede_EntreajudaLabur"
$ws2.Cells.Item(1, 2).Value = 0.3160079764148219
$ws2.Cells.Item(2, 1).Value = "OutDeg_Var.Dep_RedeControlExtAusencia"
$ws2.Cells.Item(2, 2).Value = 0.4313945633671408
$ws2.Cells.Item(3, 1).Value = "OutDeg_Var.Dep_ApoioSpecDiqCint"
$ws2.Cells.Item(3, 2).Value = 0.8525022872986581
$ws2.Cells.Item(4, 1).Value = "OutDeg_Var.Dep_ApoioSpecTubos"
$ws2.Cells.Item(4, 2).Value = 0.9760724884524125
$ws2.Cells.Item(5, 1).Value = "OutDeg_Var.Dep_RedApoiLevntDiCin"
$ws2.Cells.Item(5, 2).Value = 0.4219011615417401
$ws2.Cells.Item(6, 1).Value = "OutDeg_RedCont_DiqPriqDiqCint"
$ws2.Cells.Item(6, 2).Value = 0.8101734623123679
$ws2.Cells.Item(7, 1).Value = "OutDeg_RedeRepar_DiqPriqDiqCint"
$ws2.Cells.Item(7, 2).Value = 0.7096451735361415
$ws2.Cells.Item(8, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqCin"
$ws2.Cells.Item(8, 2).Value = 0.6035124763514501
$ws2.Cells.Item(9, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqPrq"
$ws2.Cells.Item(9, 2).Value = 0.7415294502916965
$ws2.Cells.Item(10, 1).Value = "OutDeg_ContRep2_DiqPrqCin"
$ws2.Cells.Item(10, 2).Value = 0.7772401631483951
$ws2.Cells.Item(11, 1).Value = "OutDeg_Var.Dep_RedePartilhaAgu"
$ws2.Cells.Item(11, 2).Value = 0.2489316239638737

# ---------- Sheet 'P-Values' ----------
$ws3 = $wb.Worksheets.Item("P-Values")
$ws3.Cells.Item(1, 1).Value = "Degree_Rede_EntreajudaLabur"
$ws3.Cells.Item(1, 2).Value = 0.3453225591032064
$ws3.Cells.Item(1, 3).Value = 0.4769184970427237
$ws3.Cells.Item(1, 4).Value = 0.5885359321770665
$ws3.Cells.Item(1, 5).Value = 0.05719665693973385
$ws3.Cells.Item(1, 6).Value = 0.3554578272353635
$ws3.Cells.Item(1, 7).Value = 0.008697171658797508
$ws3.Cells.Item(1, 8).Value = 0.8979779732219562
$ws3.Cells.Item(1, 9).Value = 0.5247313010925606
$ws3.Cells.Item(1, 10).Value = 0.4556506642802596
$ws3.Cells.Item(1, 11).Value = 0.7467867035722033
$ws3.Cells.Item(1, 12).Value = 0.4180288551396678
$ws3.Cells.Item(1, 13).Value = 0.1138287023618562
$ws3.Cells.Item(1, 14).Value = 0.005251602955196127
$ws3.Cells.Item(1, 15).Value = 0.1326444461424313
$ws3.Cells.Item(1, 16).Value = 0.00007545870525363093
$ws3.Cells.Item(1, 17).Value = 0.1786757775366458
$ws3.Cells.Item(1, 18).Value = 0.00000007763078758136046
$ws3.Cells.Item(1, 19).Value = 0.6760420036257873
$ws3.Cells.Item(1, 20).Value = 0.07061695485049262
$ws3.Cells.Item(1, 21).Value = 0.0002985898679668298
$ws3.Cells.Item(1, 22).Value = 0.285913068779434
$ws3.Cells.Item(1, 23).Value = 0.1135421875108657
$ws3.Cells.Item(1, 24).Value = 0.1887279861884955
$ws3.Cells.Item(1, 25).Value = 0.9207904804944991
$ws3.Cells.Item(1, 26).Value = 0.004678977644543249
$ws3.Cells.Item(1, 27).Value = 0.09702400151404798
$ws3.Cells.Item(1, 28).Value = 0.4048080317409195
$ws3.Cells.Item(1, 29).Value = 0.4573220343823444
$ws3.Cells.Item(1, 30).Value = 0.5129395155487415
$ws3.Cells.Item(2, 1).Value = "OutDeg_Var.Dep_RedeControlExtAusencia"
$ws3.Cells.Item(2, 2).Value = 0.4621986438146836
$ws3.Cells.Item(2, 3).Value = 0.659240209568646
$ws3.Cells.Item(2, 4).Value = 0.6935171864656355
$ws3.Cells.Item(2, 5).Value = 0.9873908063369523
$ws3.Cells.Item(2, 6).Value = 0.3789014969753116
$ws3.Cells.Item(2, 7).Value = 0.145035207871785
$ws3.Cells.Item(2, 8).Value = 0.001611717284525267
$ws3.Cells.Item(2, 9).Value = 0.2628267861719159
$ws3.Cells.Item(2, 10).Value = 0.8085016772695405
$ws3.Cells.Item(2, 11).Value = 0.3052987485872847
$ws3.Cells.Item(2, 12).Value = 0.8323886201457759
$ws3.Cells.Item(2, 13).Value = 0.2828944320246886
$ws3.Cells.Item(2, 14).Value = 0.8635352580058836
$ws3.Cells.Item(2, 15).Value = 0.1428662287742584
$ws3.Cells.Item(2, 16).Value = 0.3656700600566107
$ws3.Cells.Item(2, 17).Value = 0.3703797351207598
$ws3.Cells.Item(2, 18).Value = 0.5260057829194889
$ws3.Cells.Item(2, 19).Value = 0.07473725955842644
$ws3.Cells.Item(2, 20).Value = 0.1377549861730321
$ws3.Cells.Item(2, 21).Value = 0.6670819350978059
$ws3.Cells.Item(2, 22).Value = 0.8336376101058887
$ws3.Cells.Item(2, 23).Value = 0.5315737914380207
$ws3.Cells.Item(2, 24).Value = 0.08149986960820405
$ws3.Cells.Item(2, 25).Value = 0.09976442217908153
$ws3.Cells.Item(2, 26).Value = 0.09439598078347965
$ws3.Cells.Item(2, 27).Value = 0.4223664500112835
$ws3.Cells.Item(2, 28).Value = 0.0007823422723702651
$ws3.Cells.Item(2, 29).Value = 0.6206298815977549
$ws3.Cells.Item(2, 30).Value = 0.05890533387191522
$ws3.Cells.Item(3, 1).Value = "OutDeg_Var.Dep_ApoioSpecDiqCint"
$ws3.Cells.Item(3, 2).Value = 0.3826915378005611
$ws3.Cells.Item(3, 3).Value = 0.2623320064410679
$ws3.Cells.Item(3, 4).Value = 0.3082335022880889
$ws3.Cells.Item(3, 5).Value = 0.518920317362404
$ws3.Cells.Item(3, 6).Value = 0.5519679328709539
$ws3.Cells.Item(3, 7).Value = 0.5608994169763263
$ws3.Cells.Item(3, 8).Value = 0.3850294616041474
$ws3.Cells.Item(3, 9).Value = 0.3861519543747925
$ws3.Cells.Item(3, 10).Value = 0.1755048123944752
$ws3.Cells.Item(3, 11).Value = 0.2420440162236446
$ws3.Cells.Item(3, 12).Value = 0.6124123000739401
$ws3.Cells.Item(3, 13).Value = 0.01204270613352378
$ws3.Cells.Item(3, 14).Value = 0.04676895021921357
$ws3.Cells.Item(3, 15).Value = 0.9736698482453987
$ws3.Cells.Item(3, 16).Value = 0.3416668050567083
$ws3.Cells.Item(3, 17).Value = 0.386401313354519
$ws3.Cells.Item(3, 18).Value = 0.883945065707926
$ws3.Cells.Item(3, 19).Value = 0.8752203944378486
$ws3.Cells.Item(3, 20).Value = 0.07525322820963777
$ws3.Cells.Item(3, 21).Value = 0.5468245457408505
$ws3.Cells.Item(3, 22).Value = 0.4908452513078001
$ws3.Cells.Item(3, 23).Value = 0.5326033199514033
$ws3.Cells.Item(3, 24).Value = 0.464917359970742
$ws3.Cells.Item(3, 25).Value = 0.4979339052013279
$ws3.Cells.Item(3, 26).Value = 0.3219225547969327
$ws3.Cells.Item(3, 27).Value = 0.480000532418642
$ws3.Cells.Item(3, 28).Value = 0.000000000000000000008940145406718197
$ws3.Cells.Item(3, 29).Value = 0.0000000000001086673215057164
$ws3.Cells.Item(3, 30).Value = 0.1209198818121082
$ws3.Cells.Item(4, 1).Value = "OutDeg_Var.Dep_ApoioSpecTubos"
$ws3.Cells.Item(4, 2).Value = 0.08819217556465778
$ws3.Cells.Item(4, 3).Value = 0.5559920531423714
$ws3.Cells.Item(4, 4).Value = 0.2861190155136155
$ws3.Cells.Item(4, 5).Value = 0.8192558361520356
$ws3.Cells.Item(4, 6).Value = 0.4916355888641161
$ws3.Cells.Item(4, 7).Value = 0.3005679688660983
$ws3.Cells.Item(4, 8).Value = 0.6090947437661332
$ws3.Cells.Item(4, 9).Value = 0.7151896433559355
$ws3.Cells.Item(4, 10).Value = 0.7358780740664048
$ws3.Cells.Item(4, 11).Value = 0.7078300321939578
$ws3.Cells.Item(4, 12).Value = 0.9998053862510251
$ws3.Cells.Item(4, 13).Value = 0.8043845962834535
$ws3.Cells.Item(4, 14).Value = 0.03200353997164793
$ws3.Cells.Item(4, 15).Value = 0.2380040643259701
$ws3.Cells.Item(4, 16).Value = 0.1343019617429104
$ws3.Cells.Item(4, 17).Value = 0.925838689557688
$ws3.Cells.Item(4, 18).Value = 0.3180975303345709
$ws3.Cells.Item(4, 19).Value = 0.1869894602999593
$ws3.Cells.Item(4, 20).Value = 0.834141404783164
$ws3.Cells.Item(4, 21).Value = 0.7964801649823617
$ws3.Cells.Item(4, 22).Value = 0.7415264864937103
$ws3.Cells.Item(4, 23).Value = 0.00008500023489682533
$ws3.Cells.Item(4, 24).Value = 0.3049534006083039
$ws3.Cells.Item(4, 25).Value = 0.3793476821356635
$ws3.Cells.Item(4, 26).Value = 0.03269884449214878
$ws3.Cells.Item(4, 27).Value = 0.6931082216723602
$ws3.Cells.Item(4, 28).Value = 0.000006775055804526258
$ws3.Cells.Item(4, 29).Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000587956741879312
$ws3.Cells.Item(4, 30).Value = 0.7601810279743619
$ws3.Cells.Item(5, 1).Value = "OutDeg_Var.Dep_RedApoiLevntDiCin"
$ws3.Cells.Item(5, 2).Value = 0.1717144347993991
$ws3.Cells.Item(5, 3).Value = 0.5598248539860323
$ws3.Cells.Item(5, 4).Value = 0.3704830972211333
$ws3.Cells.Item(5, 5).Value = 0.9192231276448619
$ws3.Cells.Item(5, 6).Value = 0.9479812650329036
$ws3.Cells.Item(5, 7).Value = 0.720904794441289
$ws3.Cells.Item(5, 8).Value = 0.1755068029293763
$ws3.Cells.Item(5, 9).Value = 0.3455601324850698
$ws3.Cells.Item(5, 10).Value = 0.01651393069187244
$ws3.Cells.Item(5, 11).Value = 0.2858242920945619
$ws3.Cells.Item(5, 12).Value = 0.01374013118510141
$ws3.Cells.Item(5, 13).Value = 0.2977940189668127
$ws3.Cells.Item(5, 14).Value = 0.7328530455767823
$ws3.Cells.Item(5, 15).Value = 0.8067199339066021
$ws3.Cells.Item(5, 16).Value = 0.5469119128670585
$ws3.Cells.Item(5, 17).Value = 0.8703098269994356
$ws3.Cells.Item(5, 18).Value = 0.3423102613237815
$ws3.Cells.Item(5, 19).Value = 0.4858828130595879
$ws3.Cells.Item(5, 20).Value = 0.08720354581430884
$ws3.Cells.Item(5, 21).Value = 0.3048694482095105
$ws3.Cells.Item(5, 22).Value = 0.7101571422504043
$ws3.Cells.Item(5, 23).Value = 0.957854410442706
$ws3.Cells.Item(5, 24).Value = 0.4432497519000282
$ws3.Cells.Item(5, 25).Value = 0.6506607732273195
$ws3.Cells.Item(5, 26).Value = 0.03753261676410441
$ws3.Cells.Item(5, 27).Value = 0.6233059233435911
$ws3.Cells.Item(5, 28).Value = 0.00009655606067514953
$ws3.Cells.Item(5, 29).Value = 0.1809723540393573
$ws3.Cells.Item(5, 30).Value = 0.8014778010422978
$ws3.Cells.Item(6, 1).Value = "OutDeg_RedCont_DiqPriqDiqCint"
$ws3.Cells.Item(6, 2).Value = 0.8012921485595855
$ws3.Cells.Item(6, 3).Value = 0.2105351635082013
$ws3.Cells.Item(6, 4).Value = 0.2022395470514276
$ws3.Cells.Item(6, 5).Value = 0.2754164861248479
$ws3.Cells.Item(6, 6).Value = 0.7648429276378909
$ws3.Cells.Item(6, 7).Value = 0.2943296886709413
$ws3.Cells.Item(6, 8).Value = 0.9069551421424398
$ws3.Cells.Item(6, 9).Value = 0.4104365754914453
$ws3.Cells.Item(6, 10).Value = 0.1609198201624075
$ws3.Cells.Item(6, 11).Value = 0.2234908601107427
$ws3.Cells.Item(6, 12).Value = 0.7082193117213444
$ws3.Cells.Item(6, 13).Value = 0.7462666171246628
$ws3.Cells.Item(6, 14).Value = 0.09151325569103577
$ws3.Cells.Item(6, 15).Value = 0.001499557987301014
$ws3.Cells.Item(6, 16).Value = 0.3921321033147084
$ws3.Cells.Item(6, 17).Value = 0.6026277175346539
$ws3.Cells.Item(6, 18).Value = 0.08684083252966536
$ws3.Cells.Item(6, 19).Value = 0.3307767759674699
$ws3.Cells.Item(6, 20).Value = 0.5764286398066643
$ws3.Cells.Item(6, 21).Value = 0.09955322197548948
$ws3.Cells.Item(6, 22).Value = 0.8193407136721257
$ws3.Cells.Item(6, 23).Value = 0.899266541105282
$ws3.Cells.Item(6, 24).Value = 0.9334308554030533
$ws3.Cells.Item(6, 25).Value = 0.8256979763361566
$ws3.Cells.Item(6, 26).Value = 0.1341066895997916
$ws3.Cells.Item(6, 27).Value = 0.04416042817681056
$ws3.Cells.Item(6, 28).Value = 0.00000001043484047892377
$ws3.Cells.Item(6, 29).Value = 0.00000000001467626491969286
$ws3.Cells.Item(6, 30).Value = 0.3017216270949778
$ws3.Cells.Item(7, 1).Value = "OutDeg_RedeRepar_DiqPriqDiqCint"
$ws3.Cells.Item(7, 2).Value = 0.3499806365534126
$ws3.Cells.Item(7, 3).Value = 0.812012418192575
$ws3.Cells.Item(7, 4).Value = 0.7512610845546748
$ws3.Cells.Item(7, 5).Value = 0.6730221961397839
$ws3.Cells.Item(7, 6).Value = 0.8207519628603337
$ws3.Cells.Item(7, 7).Value = 0.1443013596364536
$ws3.Cells.Item(7, 8).Value = 0.2259768653205903
$ws3.Cells.Item(7, 9).Value = 0.5628963125129953
$ws3.Cells.Item(7, 10).Value = 0.8568504647113053
$ws3.Cells.Item(7, 11).Value = 0.8880682819482866
$ws3.Cells.Item(7, 12).Value = 0.7888359943679206
$ws3.Cells.Item(7, 13).Value = 0.6759752966387632
$ws3.Cells.Item(7, 14).Value = 0.07221663189387137
$ws3.Cells.Item(7, 15).Value = 0.06604590957384113
$ws3.Cells.Item(7, 16).Value = 0.6755525536254218
$ws3.Cells.Item(7, 17).Value = 0.4196488004828974
$ws3.Cells.Item(7, 18).Value = 0.1073186237273558
$ws3.Cells.Item(7, 19).Value = 0.1978187960291865
$ws3.Cells.Item(7, 20).Value = 0.3956018932906783
$ws3.Cells.Item(7, 21).Value = 0.3626801617572745
$ws3.Cells.Item(7, 22).Value = 0.6551681221715535
$ws3.Cells.Item(7, 23).Value = 0.8789449218694814
$ws3.Cells.Item(7, 24).Value = 0.6184375745381638
$ws3.Cells.Item(7, 25).Value = 0.6460294348107135
$ws3.Cells.Item(7, 26).Value = 0.005078748675022629
$ws3.Cells.Item(7, 27).Value = 0.2646768132302947
$ws3.Cells.Item(7, 28).Value = 0.00000001757894937209405
$ws3.Cells.Item(7, 29).Value = 0.000001320351234845503
$ws3.Cells.Item(7, 30).Value = 0.1441737534085736
$ws3.Cells.Item(8, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqCin"
$ws3.Cells.Item(8, 2).Value = 0.4108443401785097
$ws3.Cells.Item(8, 3).Value = 0.7809128771075595
$ws3.Cells.Item(8, 4).Value = 0.6561623548179951
$ws3.Cells.Item(8, 5).Value = 0.6778380950880042
$ws3.Cells.Item(8, 6).Value = 0.453081346070272
$ws3.Cells.Item(8, 7).Value = 0.838759514015363
$ws3.Cells.Item(8, 8).Value = 0.581096360177596
$ws3.Cells.Item(8, 9).Value = 0.7222601518516061
$ws3.Cells.Item(8, 10).Value = 0.09523353910384029
$ws3.Cells.Item(8, 11).Value = 0.80430636756558
$ws3.Cells.Item(8, 12).Value = 0.789995456026612
$ws3.Cells.Item(8, 13).Value = 0.4747895358497193
$ws3.Cells.Item(8, 14).Value = 0.09928515891087301
$ws3.Cells.Item(8, 15).Value = 0.02713640161951018
$ws3.Cells.Item(8, 16).Value = 0.7200046238451665
$ws3.Cells.Item(8, 17).Value = 0.01063683491469571
$ws3.Cells.Item(8, 18).Value = 0.1250646551320441
$ws3.Cells.Item(8, 19).Value = 0.347876112368069
$ws3.Cells.Item(8, 20).Value = 0.1294451630131721
$ws3.Cells.Item(8, 21).Value = 0.3867599548731451
$ws3.Cells.Item(8, 22).Value = 0.6227808488765427
$ws3.Cells.Item(8, 23).Value = 0.1748946341311294
$ws3.Cells.Item(8, 24).Value = 0.99573950800826
$ws3.Cells.Item(8, 25).Value = 0.7282439810910786
$ws3.Cells.Item(8, 26).Value = 0.1336173398351177
$ws3.Cells.Item(8, 27).Value = 0.08801738596310298
$ws3.Cells.Item(8, 28).Value = 0.0000003330157462800386
$ws3.Cells.Item(8, 29).Value = 0.5009190543867861
$ws3.Cells.Item(8, 30).Value = 0.5041565843744081
$ws3.Cells.Item(9, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqPrq"
$ws3.Cells.Item(9, 2).Value = 0.5313045714230722
$ws3.Cells.Item(9, 3).Value = 0.4062907811069032
$ws3.Cells.Item(9, 4).Value = 0.3105474100195874
$ws3.Cells.Item(9, 5).Value = 0.2358305080711429
$ws3.Cells.Item(9, 6).Value = 0.793641693906628
$ws3.Cells.Item(9, 7).Value = 0.138152575903449
$ws3.Cells.Item(9, 8).Value = 0.3350788909689879
$ws3.Cells.Item(9, 9).Value = 0.7881141683120586
$ws3.Cells.Item(9, 10).Value = 0.6249207880313077
$ws3.Cells.Item(9, 11).Value = 0.270733734734854
$ws3.Cells.Item(9, 12).Value = 0.6943384066690244
$ws3.Cells.Item(9, 13).Value = 0.9461497936979693
$ws3.Cells.Item(9, 14).Value = 0.08309376796766403
$ws3.Cells.Item(9, 15).Value = 0.01746530193194068
$ws3.Cells.Item(9, 16).Value = 0.5612012880687745
$ws3.Cells.Item(9, 17).Value = 0.7560673224350793
$ws3.Cells.Item(9, 18).Value = 0.1472930517313237
$ws3.Cells.Item(9, 19).Value = 0.1561926367840621
$ws3.Cells.Item(9, 20).Value = 0.4452840810108318
$ws3.Cells.Item(9, 21).Value = 0.2855223607641528
$ws3.Cells.Item(9, 22).Value = 0.4471148900120298
$ws3.Cells.Item(9, 23).Value = 0.5527950094567649
$ws3.Cells.Item(9, 24).Value = 0.6080012432683596
$ws3.Cells.Item(9, 25).Value = 0.6624593225652056
$ws3.Cells.Item(9, 26).Value = 0.07576059571746975
$ws3.Cells.Item(9, 27).Value = 0.18421963259142
$ws3.Cells.Item(9, 28).Value = 0.0000008456122686816655
$ws3.Cells.Item(9, 29).Value = 0.000000001382234064552799
$ws3.Cells.Item(9, 30).Value = 0.1872250389785442
$ws3.Cells.Item(10, 1).Value = "OutDeg_ContRep2_DiqPrqCin"
$ws3.Cells.Item(10, 2).Value = 0.5109679341761201
$ws3.Cells.Item(10, 3).Value = 0.5097383333604346
$ws3.Cells.Item(10, 4).Value = 0.4542089533359458
$ws3.Cells.Item(10, 5).Value = 0.5064431738824555
$ws3.Cells.Item(10, 6).Value = 0.6404805860977896
$ws3.Cells.Item(10, 7).Value = 0.1640354416301944
$ws3.Cells.Item(10, 8).Value = 0.2923519649402742
$ws3.Cells.Item(10, 9).Value = 0.7974054577939003
$ws3.Cells.Item(10, 10).Value = 0.5245543505300501
$ws3.Cells.Item(10, 11).Value = 0.7175999108193681
$ws3.Cells.Item(10, 12).Value = 0.9076258844966294
$ws3.Cells.Item(10, 13).Value = 0.6264446170430613
$ws3.Cells.Item(10, 14).Value = 0.0660866172445911
$ws3.Cells.Item(10, 15).Value = 0.0225410477181911
$ws3.Cells.Item(10, 16).Value = 0.4853937602773013
$ws3.Cells.Item(10, 17).Value = 0.4824835263869158
$ws3.Cells.Item(10, 18).Value = 0.06595346612673328
$ws3.Cells.Item(10, 19).Value = 0.2208032343566827
$ws3.Cells.Item(10, 20).Value = 0.4169669725522683
$ws3.Cells.Item(10, 21).Value = 0.1937185136142919
$ws3.Cells.Item(10, 22).Value = 0.7789384630905081
$ws3.Cells.Item(10, 23).Value = 0.9659470983530244
$ws3.Cells.Item(10, 24).Value = 0.6729124819559584
$ws3.Cells.Item(10, 25).Value = 0.7362573980032031
$ws3.Cells.Item(10, 26).Value = 0.01587914306795081
$ws3.Cells.Item(10, 27).Value = 0.09231230483432042
$ws3.Cells.Item(10, 28).Value = 0.000000004360023492404507
$ws3.Cells.Item(10, 29).Value = 0.000000002375011066748857
$ws3.Cells.Item(10, 30).Value = 0.1725523806655892
$ws3.Cells.Item(11, 1).Value = "OutDeg_Var.Dep_RedePartilhaAgu"
$ws3.Cells.Item(11, 2).Value = 0.2834711344987654
$ws3.Cells.Item(11, 3).Value = 0.0584218437406451
$ws3.Cells.Item(11, 4).Value = 0.2139088873037175
$ws3.Cells.Item(11, 5).Value = 0.4662707961653342
$ws3.Cells.Item(11, 6).Value = 0.1916051891340759
$ws3.Cells.Item(11, 7).Value = 0.2948313847455645
$ws3.Cells.Item(11, 8).Value = 0.5318451579701762
$ws3.Cells.Item(11, 9).Value = 0.4585021984655258
$ws3.Cells.Item(11, 10).Value = 0.7090718178182777
$ws3.Cells.Item(11, 11).Value = 0.004129310139171982
$ws3.Cells.Item(11, 12).Value = 0.07292380758530818
$ws3.Cells.Item(11, 13).Value = 0.1635625119476668
$ws3.Cells.Item(11, 14).Value = 0.1469659180518643
$ws3.Cells.Item(11, 15).Value = 0.2311294662266287
$ws3.Cells.Item(11, 16).Value = 0.875503578246148
$ws3.Cells.Item(11, 17).Value = 0.4891976973132139
$ws3.Cells.Item(11, 18).Value = 0.8820749063656508
$ws3.Cells.Item(11, 19).Value = 0.871215752281697
$ws3.Cells.Item(11, 20).Value = 0.9960260365738505
$ws3.Cells.Item(11, 21).Value = 0.1260550011328726
$ws3.Cells.Item(11, 22).Value = 0.755904719899004
$ws3.Cells.Item(11, 23).Value = 0.7856895361670597
$ws3.Cells.Item(11, 24).Value = 0.6464895972583977
$ws3.Cells.Item(11, 25).Value = 0.7045032883327984
$ws3.Cells.Item(11, 26).Value = 0.5468699763527667
$ws3.Cells.Item(11, 27).Value = 0.3265615861276173
$ws3.Cells.Item(11, 28).Value = 0.3996426637955943
$ws3.Cells.Item(11, 29).Value = 0.1508235780084675
$ws3.Cells.Item(11, 30).Value = 0.9220651405294855

# ---------- Sheet 'Regresiones con R2 > 0.1' ----------
$ws4 = $wb.Worksheets.Item("Regresiones con R2 > 0.1")
$ws4.Cells.Item(1, 1).Value = "Degree_Rede_EntreajudaLabur"
$ws4.Cells.Item(1, 2).Value = 0.3160079764148219
$ws4.Cells.Item(2, 1).Value = "OutDeg_Var.Dep_RedeControlExtAusencia"
$ws4.Cells.Item(2, 2).Value = 0.4313945633671408
$ws4.Cells.Item(3, 1).Value = "OutDeg_Var.Dep_ApoioSpecDiqCint"
$ws4.Cells.Item(3, 2).Value = 0.8525022872986581
$ws4.Cells.Item(4, 1).Value = "OutDeg_Var.Dep_ApoioSpecTubos"
$ws4.Cells.Item(4, 2).Value = 0.9760724884524125
$ws4.Cells.Item(5, 1).Value = "OutDeg_Var.Dep_RedApoiLevntDiCin"
$ws4.Cells.Item(5, 2).Value = 0.4219011615417401
$ws4.Cells.Item(6, 1).Value = "OutDeg_RedCont_DiqPriqDiqCint"
$ws4.Cells.Item(6, 2).Value = 0.8101734623123679
$ws4.Cells.Item(7, 1).Value = "OutDeg_RedeRepar_DiqPriqDiqCint"
$ws4.Cells.Item(7, 2).Value = 0.7096451735361415
$ws4.Cells.Item(8, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqCin"
$ws4.Cells.Item(8, 2).Value = 0.6035124763514501
$ws4.Cells.Item(9, 1).Value = "OutDeg_Red-Val_Cont1Rep2DiqPrq"
$ws4.Cells.Item(9, 2).Value = 0.7415294502916965
$ws4.Cells.Item(10, 1).Value = "OutDeg_ContRep2_DiqPrqCin"
$ws4.Cells.Item(10, 2).Value = 0.7772401631483951
$ws4.Cells.Item(11, 1).Value = "OutDeg_Var.Dep_RedePartilhaAgu"
$ws4.Cells.Item(11, 2).Value = 0.2489316239638737

Write-Output "done"